$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2 and H3 need scientific-notation precision; compute via formula then bake to a static value.
$ws.Range("H2").Formula = "=1.41429684665625*10^-16"
$h2 = $ws.Range("H2").Value2
$ws.Range("H2").Value = $h2

$ws.Range("H3").Formula = "=1.41429684665625*10^-16"
$h3 = $ws.Range("H3").Value2
$ws.Range("H3").Value = $h3

# Row 2
$ws.Range("E2").Value = 25.93000000000061
$ws.Range("K2").Value = 48.19965900304013
$ws.Range("L2").Value = "[45.606676500374164, 50.79264150570609]"
$ws.Range("O2").Value = 1.603816069400195
$ws.Range("P2").Value = "[1.553500271144502, 1.654131867655888]"
$ws.Range("S2").Value = 52.51806210883431
$ws.Range("T2").Value = "[50.758889139715926, 54.27723507795269]"
$ws.Range("W2").Value = 19.31123123123169
$ws.Range("X2").Value = 19.10358358358403
$ws.Range("Y2").Value = 19.51887887887934

# Row 3
$ws.Range("E3").Value = 22.79000000000012
$ws.Range("K3").Value = 48.65731167234735
$ws.Range("L3").Value = "[43.098198848418726, 54.216424496275984]"
$ws.Range("O3").Value = 2.207605648468503
$ws.Range("P3").Value = "[2.094395102393195, 2.320816194543811]"
$ws.Range("S3").Value = 52.16440885480622
$ws.Range("T3").Value = "[49.321443198291924, 55.00737451132051]"
$ws.Range("W3").Value = 14.78270270270278
$ws.Range("X3").Value = 14.37207207207215
$ws.Range("Y3").Value = 15.19333333333342
